$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2230576441102757
$ws.Range("C2").Value = 0.4987468671679198
$ws.Range("J2").Value = 0.012531328320802
$ws.Range("P2").Value = 0.1528822055137845
$ws.Range("S2").Value = 0.112781954887218
$ws.Range("B3").Value = 0.009852216748768473
$ws.Range("C3").Value = 0.01477832512315271
$ws.Range("J3").Value = 0.01477832512315271
$ws.Range("P3").Value = 0.729064039408867
$ws.Range("S3").Value = 0.2315270935960591
$ws.Range("J4").Value = 0.05357142857142857
$ws.Range("P4").Value = 0.5892857142857143
$ws.Range("S4").Value = 0.3571428571428572
$ws.Range("B6").Value = 0.07865168539325842
$ws.Range("D6").Value = 0.02247191011235955
$ws.Range("F6").Value = 0.09363295880149813
$ws.Range("J6").Value = 0.2846441947565543
$ws.Range("O6").Value = 0.0149812734082397
$ws.Range("Q6").Value = 0.1460674157303371
$ws.Range("R6").Value = 0.08239700374531835
$ws.Range("S6").Value = 0.2771535580524345
$ws.Range("B7").Value = 0.1176470588235294
$ws.Range("D7").Value = 0.02941176470588235
$ws.Range("E7").Value = 0.007352941176470588
$ws.Range("F7").Value = 0.04411764705882353
$ws.Range("J7").Value = 0.1580882352941176
$ws.Range("O7").Value = 0.01470588235294118
$ws.Range("Q7").Value = 0.1985294117647059
$ws.Range("R7").Value = 0.1102941176470588
$ws.Range("S7").Value = 0.3198529411764706
$ws.Range("B8").Value = 0.09861932938856016
$ws.Range("D8").Value = 0.01577909270216963
$ws.Range("F8").Value = 0.0611439842209073
$ws.Range("J8").Value = 0.1203155818540434
$ws.Range("O8").Value = 0.01577909270216963
$ws.Range("Q8").Value = 0.1577909270216963
$ws.Range("R8").Value = 0.1380670611439842
$ws.Range("S8").Value = 0.3925049309664694
$ws.Range("B9").Value = 0.1016949152542373
$ws.Range("D9").Value = 0.0211864406779661
$ws.Range("F9").Value = 0.05932203389830509
$ws.Range("J9").Value = 0.09322033898305085
$ws.Range("O9").Value = 0.01694915254237288
$ws.Range("Q9").Value = 0.1186440677966102
$ws.Range("R9").Value = 0.1186440677966102
$ws.Range("S9").Value = 0.4703389830508475
$ws.Range("B10").Value = 0.1171112556929083
$ws.Range("D10").Value = 0.01951854261548471
$ws.Range("F10").Value = 0.06115810019518542
$ws.Range("J10").Value = 0.1301236174365647
$ws.Range("O10").Value = 0.01431359791802212
$ws.Range("Q10").Value = 0.214053350683149
$ws.Range("R10").Value = 0.08327911515940144
$ws.Range("S10").Value = 0.3604424202992843
$ws.Range("G11").Value = 0.1471264367816092
$ws.Range("J11").Value = 0.103448275862069
$ws.Range("K11").Value = 0.2068965517241379
$ws.Range("L11").Value = 0.535632183908046
$ws.Range("S11").Value = 0.006896551724137931
$ws.Range("G12").Value = 0.7269076305220884
$ws.Range("J12").Value = 0.1927710843373494
$ws.Range("K12").Value = 0.01204819277108434
$ws.Range("L12").Value = 0.04417670682730924
$ws.Range("S12").Value = 0.02409638554216868
$ws.Range("G13").Value = 0.660377358490566
$ws.Range("J13").Value = 0.3018867924528302
$ws.Range("S13").Value = 0.03773584905660377
$ws.Range("F15").Value = 0.03515625
$ws.Range("H15").Value = 0.171875
$ws.Range("I15").Value = 0.0703125
$ws.Range("J15").Value = 0.32421875
$ws.Range("K15").Value = 0.08203125
$ws.Range("M15").Value = 0.01171875
$ws.Range("O15").Value = 0.11328125
$ws.Range("S15").Value = 0.19140625
$ws.Range("F16").Value = 0.01304347826086956
$ws.Range("H16").Value = 0.1521739130434783
$ws.Range("I16").Value = 0.08695652173913043
$ws.Range("J16").Value = 0.4434782608695652
$ws.Range("K16").Value = 0.1304347826086956
$ws.Range("M16").Value = 0.008695652173913044
$ws.Range("O16").Value = 0.03478260869565217
$ws.Range("S16").Value = 0.1304347826086956
$ws.Range("F17").Value = 0.01318267419962335
$ws.Range("H17").Value = 0.1431261770244821
$ws.Range("I17").Value = 0.09981167608286252
$ws.Range("J17").Value = 0.384180790960452
$ws.Range("K17").Value = 0.1337099811676083
$ws.Range("M17").Value = 0.032015065913371
$ws.Range("N17").Value = 0.001883239171374765
$ws.Range("O17").Value = 0.064030131826742
$ws.Range("S17").Value = 0.128060263653484
$ws.Range("F18").Value = 0.03249097472924187
$ws.Range("H18").Value = 0.1732851985559567
$ws.Range("I18").Value = 0.08664259927797834
$ws.Range("J18").Value = 0.3898916967509025
$ws.Range("K18").Value = 0.1083032490974729
$ws.Range("M18").Value = 0.01805054151624549
$ws.Range("O18").Value = 0.05054151624548736
$ws.Range("S18").Value = 0.1407942238267148
$ws.Range("F19").Value = 0.02220726783310902
$ws.Range("H19").Value = 0.2072678331090175
$ws.Range("I19").Value = 0.08344549125168237
$ws.Range("J19").Value = 0.360699865410498
$ws.Range("K19").Value = 0.1244952893674293
$ws.Range("M19").Value = 0.01816958277254374
$ws.Range("N19").Value = 0.0006729475100942127
$ws.Range("O19").Value = 0.06729475100942127
$ws.Range("S19").Value = 0.1157469717362046
